$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMP_ATTENDANCE")

# New attendance row appended below the existing data (A1:C5 -> A1:C6).
# Source values look numeric/date-like ("32", "05-07-2022") but the sheet
# stores everything as plain text, so force text formatting before writing
# them to avoid Excel auto-converting to a number/date serial.
$ws.Range("A6:B6").NumberFormat = "@"

$ws.Range("A6").Value = "32"
$ws.Range("B6").Value = "05-07-2022"
$ws.Range("C6").Value = "20:13:08"
